# "ultima version de los gimnasios" - update the Cucuta gyms listing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "SUPER GIMNASIO LA 13" is replaced by "BODY TECH"
# (write D3 first, then B3, then C3 so the shared-string table grows
#  in the same order the source workbook used)
$ws.Range("D3").Value = "7.887723, -72.494997"
$ws.Range("B3").Value = "BODY TECH"
$ws.Range("C3").Value = "Cl. 11 #2E-10 caobos"

# Row 4: CACIQUE FITNNES address re-cased
$ws.Range("C4").Value = "AV 2 16-01 la playa"

# Row 10: ZONA DE PILATES address updated
$ws.Range("C10").Value = "Av. 1 Este #2022 Barrio blanco"

# Row 13: ZONA NORTE MMA CUCUTA address typo fixed (## -> #)
$ws.Range("C13").Value = "Av. 3 Este # 14a, Cúcuta, Norte de Santander"

# Row 14: GYM LA 10 address updated
$ws.Range("C14").Value = "8- a, Av. 10 #8125, barrio llano"

# Restore selection to C18, matching the saved workbook view
$ws.Range("C18").Select()
